$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table
$tbl.Rows.Add()
$newRow = $tbl.Rows.Count
$tbl.Cell($newRow, 1).Shape.TextFrame.TextRange.Text = "Car"
$tbl.Cell($newRow, 2).Shape.TextFrame.TextRange.Text = "75"
$tbl.Cell($newRow, 3).Shape.TextFrame.TextRange.Text = "2025-01-01"
$tbl.Cell($newRow, 4).Shape.TextFrame.TextRange.Text = "Successfully added expense"
